# Revise build script optimization
# Appends newly-collected sensor samples to the two sensor pairs
# (ROW50 / ROW11, each FE-LIFTER + MID-LIFTER) and converts the prior
# "latest reading" placeholder row (which was stored as plain text)
# into a proper numeric/date-formatted row now that it has settled.

$wb = $excel.ActiveWorkbook

$DATE_FMT = "YYYY-MM-DD HH:MM:SS"
$G_VALUE = [double]"5.686312626471138e+23"

function Set-SensorRow($ws, $r, $a, $c, $e, $i) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 1).NumberFormat = $DATE_FMT
    $ws.Cells.Item($r, 2).Value = "0x01,0x90"
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = "0x01,0x90,"
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = 400
    $ws.Cells.Item($r, 7).Value = $G_VALUE
    $ws.Cells.Item($r, 8).Value = 400
    $ws.Cells.Item($r, 9).Value = $i
}

function Set-SensorRowText($ws, $r, $aText, $c, $e, $i) {
    $ws.Cells.Item($r, 1).Value = $aText
    $ws.Cells.Item($r, 2).Value = "0x01,0x90"
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = "0x01,0x90,"
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = 400
    $ws.Cells.Item($r, 7).Value = $G_VALUE
    $ws.Cells.Item($r, 8).Value = 400
    $ws.Cells.Item($r, 9).Value = $i
}

$feId = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$feChk = "0x14"
$midId = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$midChk = "0x19"

# --- ROW50-FE-LIFTER (sheet 1) and ROW11-FE-LIFTER (sheet 3) ---
# Both sheets share the identical edit: row 26 (the trailing text
# placeholder) becomes a real numeric/date value, two more settled
# numeric rows (27,28) are appended, and a fresh trailing text
# placeholder row (29) is added for the newest reading.
$feSheets = @(1, 3)
foreach ($idx in $feSheets) {
    $ws = $wb.Worksheets.Item($idx)

    Set-SensorRow $ws 26 45729.08020857639 $feId $feChk 20
    Set-SensorRow $ws 27 45729.08023003472 $feId $feChk 20
    Set-SensorRow $ws 28 45729.08025329861 $feId $feChk 20
    Set-SensorRowText $ws 29 "2025-03-13 13:55:33" $feId $feChk 20
}

# --- ROW50-MID-LIFTER (sheet 2) and ROW11-MID-LIFTER (sheet 4) ---
# Both sheets share the identical edit: nine new settled numeric rows
# (59-67) are appended after the existing data.
$midSheets = @(2, 4)
$midTimes = @(
    45729.06475774306,
    45729.06477990741,
    45729.06480305555,
    45729.14823328704,
    45729.14825528935,
    45729.14827864584,
    45729.23170871528,
    45729.23173081018,
    45729.23175395833
)

foreach ($idx in $midSheets) {
    $ws = $wb.Worksheets.Item($idx)

    $r = 59
    foreach ($t in $midTimes) {
        Set-SensorRow $ws $r $t $midId $midChk 25
        $r = $r + 1
    }
}
